# Update the "想去人数" (people interested) counts in column F
# for both the "展览" and "全部类型" worksheets, which contain the
# same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 2641
    $ws.Range("F4").Value = 536
    $ws.Range("F6").Value = 6603
    $ws.Range("F7").Value = 453
    $ws.Range("F10").Value = 12
    $ws.Range("F11").Value = 5
}
